$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header row (shelter-to-barangay distance matrix labels)
$ws.Range("A1").Value = "Shelters"
$ws.Range("B1").Value = "Balite"
$ws.Range("C1").Value = "Balungao"
$ws.Range("D1").Value = "Bulusan"
$ws.Range("E1").Value = "Frances"
$ws.Range("F1").Value = "Gatbuca"
$ws.Range("G1").Value = "Iba O'Este"

# Row 2 - Gatbuca Basketball Court (values updated)
$ws.Range("A2").Value = "Gatbuca Basketball Court"
$ws.Range("B2").Value = 4.777907
$ws.Range("C2").Value = 1.46901
$ws.Range("D2").Value = 3.221099999999999
$ws.Range("E2").Value = 1.913019
$ws.Range("F2").Value = 0.695246
$ws.Range("G2").Value = 3.994042999999999

# Row 3 - San Miguel Meysulao High School
$ws.Range("A3").Value = "San Miguel Meysulao High School"
$ws.Range("B3").Value = 7.226274999999998
$ws.Range("C3").Value = 3.917378
$ws.Range("D3").Value = 5.669467999999999
$ws.Range("E3").Value = 1.262421
$ws.Range("F3").Value = 3.216882
$ws.Range("G3").Value = 6.442410999999998

# Row 4 - Doña Damiana Elem School
$ws.Range("A4").Value = "Doña Damiana Elem School"
$ws.Range("B4").Value = 7.226274999999998
$ws.Range("C4").Value = 3.917378
$ws.Range("D4").Value = 5.669467999999999
$ws.Range("E4").Value = 1.262421
$ws.Range("F4").Value = 3.216882
$ws.Range("G4").Value = 6.442410999999998

# Row 5 - Calizon Dike (new row)
$ws.Range("A5").Value = "Calizon Dike"
$ws.Range("B5").Value = 4.782519000000001
$ws.Range("C5").Value = 0.67944
$ws.Range("D5").Value = 1.381525
$ws.Range("E5").Value = 3.104334
$ws.Range("F5").Value = 1.959829
$ws.Range("G5").Value = 3.998654999999999

# Row 6 - Frances E.C. (new row)
$ws.Range("A6").Value = "Frances E.C."
$ws.Range("B6").Value = 4.778893999999999
$ws.Range("C6").Value = 1.469997
$ws.Range("D6").Value = 3.222087
$ws.Range("E6").Value = 1.18496
$ws.Range("F6").Value = 0.769501
$ws.Range("G6").Value = 3.99503

# Apply the existing header-cell formatting (bold, thin border, centered) to the
# newly added header cells F1:G1, matching the rest of row 1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Apply the existing shelter-label formatting (bold, thin border, centered) to
# the newly added rows' column-A cells, matching A2:A4
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
